# KIBON-1710 increase erlaeuterungZ10 key, since 0 seems not to work with POI
#
# The sheet has a small "Erlaeuterungen" (explanations) legend at the bottom
# (rows 15-22). The row that was labelled with the placeholder key
# "erlaeuterungZ1_*" actually documents explanation line "12", but because
# its key was just "erlaeuterungZ1" it collided (as a prefix) with the
# differently-used "erlaeuterungZ1" placeholder elsewhere, which broke POI's
# substitution. The row is renumbered/renamed to "erlaeuterungZ12_*" and
# moved after the "erlaeuterungZ11_*" row so the rows stay in numeric order.
# Along with this, the legend block gets a lighter/smaller (9pt) font and a
# bold heading row, and the explanatory row (18) becomes a bit taller to fit
# wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the content of row 21 and row 22, and rename the (formerly row 21)
#    "erlaeuterungZ1_1"/"erlaeuterungZ1_2" placeholders to
#    "erlaeuterungZ12_1"/"erlaeuterungZ12_2".
# ---------------------------------------------------------------------------
$oldRow21A = $ws.Range("A21").Value()
$oldRow21B = $ws.Range("B21").Value()
$oldRow22A = $ws.Range("A22").Value()
$oldRow22B = $ws.Range("B22").Value()

$renamedRow21A = $oldRow21A.Replace("{erlaeuterungZ1_1}", "{erlaeuterungZ12_1}")
$renamedRow21B = $oldRow21B.Replace("{erlaeuterungZ1_2}", "{erlaeuterungZ12_2}")

$ws.Range("A21").Value = $oldRow22A
$ws.Range("B21").Value = $oldRow22B
$ws.Range("A22").Value = $renamedRow21A
$ws.Range("B22").Value = $renamedRow21B

# ---------------------------------------------------------------------------
# 2) Make the three parameter rows above the table (A3:A5) use the small
#    9pt font.
# ---------------------------------------------------------------------------
$ws.Range("A3:A5").Font.Size = 9

# ---------------------------------------------------------------------------
# 3) New legend heading row (row 15) gets formatted cells: bold 9pt for the
#    data columns, with column C right-aligned.
# ---------------------------------------------------------------------------
$ws.Range("A15").Font.Size = 11

$ws.Range("B15:J15").Font.Size = 9
$ws.Range("B15:J15").Font.Bold = $true

$ws.Range("C15").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 4) Legend rows 16-22 (the explanation table) move from 11pt (bold for
#    column A only) to a consistent bold 9pt font across all their cells.
# ---------------------------------------------------------------------------
$ws.Range("A16:A22").Font.Size = 9
$ws.Range("A16:A22").Font.Bold = $true

$ws.Range("B16:J22").Font.Size = 9
$ws.Range("B16:J22").Font.Bold = $true

# ---------------------------------------------------------------------------
# 5) Row 18's explanation text needs more vertical room once wrapped at the
#    smaller font.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).RowHeight = 29.25

# ---------------------------------------------------------------------------
# 6) Leave the cursor at A1 (the previous save had a stale selection sitting
#    on the old A21).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
